# Active_Outages.xlsx update — 6/19/2025, 8:19:08 AM
# 1) Remove the stale "HAJ0155" outage row from the R1 sheet.
# 2) Refresh the "Elapsed Duration(Hrs)" values across the remaining open
#    outages to reflect the later report-generation time.

$wb = $excel.ActiveWorkbook

# --- 1. Delete row 6 on sheet "R1" (PCM blank / HAJ0155 / Good) ---
$wsR1 = $wb.Worksheets.Item("R1")
$wsR1.Rows.Item(6).Delete()

# --- 2. Update elapsed-duration (column G) values ---
$wsR1.Range("G2").Value = "3945:33:21"
$wsR1.Range("G3").Value = "85:05:59"
$wsR1.Range("G4").Value = "108:05:59"

$wsR2 = $wb.Worksheets.Item("R2")
$wsR2.Range("G2").Value = "12126:57:02"
$wsR2.Range("G3").Value = "3256:40:31"
$wsR2.Range("G4").Value = "494:52:05"

$wsR4 = $wb.Worksheets.Item("R4")
$wsR4.Range("G2").Value = "2972:46:51"
$wsR4.Range("G3").Value = "199:59:06"
$wsR4.Range("G4").Value = "88:11:31"
$wsR4.Range("G5").Value = "85:49:04"

$wsR5 = $wb.Worksheets.Item("R5")
$wsR5.Range("G2").Value = "446:45:50"

$wsR6 = $wb.Worksheets.Item("R6")
$wsR6.Range("G2").Value = "87:18:08"
